# Daily attendance processing - reorder "Recorded By" (column G) author lists.
#
# For a known set of author-list strings, the order of the comma separated
# entries is swapped around (e.g. "System, dnasr281@gmail.com" becomes
# "dnasr281@gmail.com, System"). We detect the current value of every cell
# in column G and, whenever it matches one of the known strings, replace it
# with its reordered counterpart. Anything else (single-author cells, the
# already-reordered "dnasr281@gmail.com, admin@admin.com" cells, headers,
# etc.) is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System"
    "backup@backdoor.com, System"         = "System, backup@backdoor.com"
    "backup@backdoor.com, System, system" = "system, System, backup@backdoor.com"
    "System, admin@admin.com"             = "admin@admin.com, System"
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 2 }

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $current = $cell.Value()
    if ($null -ne $current -and $map.ContainsKey($current)) {
        $cell.Value = $map[$current]
    }
}
